# "added 250 packet run"
#
# The workbook already had a "250 Packets" column group (M:O) with its
# Acc/Loss/Time sub-headers in place but no trial data yet. This fills in
# the 49 trial rows (rows 3-51) for that run. It also starts a new,
# header-only "230 Packets" group (Q:S) — mirroring the A:C/E:G/I:K/M:O
# layout — for a future run whose data hasn't been entered yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "230 Packets" group header (row 1) and Acc/Loss/Time sub-headers (row 2)
$ws.Range("R1").Value = "230 Packets"
$ws.Range("Q2").Value = "Acc"
$ws.Range("R2").Value = "Loss"
$ws.Range("S2").Value = "Time"

# 49 trial rows (rows 3-51) of Accuracy / Loss / Time for the "250 Packets"
# run, written into columns M, N, O respectively.
$accVals  = @(89.1787409782409, 89.069652557373, 89.3660902976989, 89.311546087265, 88.6143207550048, 89.2522573471069, 89.444351196289, 89.1763687133789, 88.713926076889, 89.2072021961212, 88.5479152202606, 89.0957355499267, 88.8372421264648, 88.3392214775085, 89.2404019832611, 88.4317100048065, 89.4846677780151, 89.107596874237, 89.6388173103332, 88.7305259704589, 88.6048316955566, 89.7479057312011, 89.5249843597412, 88.7281537055969, 89.3732070922851, 88.6546373367309, 88.8894140720367, 89.7668778896331, 90.1676654815673, 89.1479134559631, 88.4981155395507, 88.9060139656066, 89.2498850822448, 89.2261743545532, 89.3732070922851, 88.939219713211, 88.974791765213, 89.3376350402832, 89.3708348274231, 88.8965308666229, 89.0411913394928, 88.0854725837707, 88.9297306537628, 89.0577912330627, 89.520239830017, 88.8372421264648, 89.1099691390991, 89.0933692455291, 89.3423795700073)
$lossVals = @(0.252545211645541, 0.252583463879403, 0.253598626512814, 0.257282045088646, 0.259854596073301, 0.250105205475315, 0.252062610138496, 0.259844600847236, 0.262258270173076, 0.261639525922112, 0.257852350244432, 0.257739163145723, 0.25607624672088, 0.254073648420198, 0.255613107552373, 0.254592751260529, 0.242295201575834, 0.253417875005239, 0.242332365250134, 0.251929461286537, 0.264360499260949, 0.24085296937197, 0.249938594938422, 0.268534642831454, 0.245401347847126, 0.25298739580283, 0.245978825268534, 0.25262072788502, 0.243155730558142, 0.260568303093277, 0.26244507051403, 0.257965461422029, 0.246497458720698, 0.249584855786408, 0.249828237407385, 0.254716357877023, 0.25970559938143, 0.250110044411053, 0.254122087232139, 0.258900934794837, 0.2572344521875, 0.260462202648749, 0.255603445054533, 0.255339804829434, 0.25332533133852, 0.260961921851552, 0.247155980114638, 0.25077409204323, 0.248062007436813)
$timeVals = @(57.434668302536, 50.7363424301147, 52.4172804355621, 70.2349917888641, 72.4229145050048, 73.7809717655181, 73.4703638553619, 73.8881390094757, 73.66739153862, 73.16757106781, 73.0641400814056, 73.3716349601745, 73.2991461753845, 72.0548400878906, 66.8185873031616, 65.3156206607818, 65.3256409168243, 64.942824602127, 71.6406962871551, 72.8096537590026, 73.1050026416778, 72.4225943088531, 73.3759329319, 72.4356048107147, 71.9085404872894, 72.832605600357, 73.3712086677551, 73.3398666381836, 72.9883384704589, 73.1948158740997, 72.4619164466857, 72.1333606243133, 72.210428237915, 73.164016008377, 72.7936882972717, 72.7490456104278, 72.1685619354248, 73.4771218299865, 74.1282415390014, 71.630872964859, 73.0956282615661, 73.0306665897369, 73.4109938144683, 72.2200510501861, 73.4147531986236, 74.0242829322815, 73.7614839076995, 72.682822227478, 73.1916553974151)

for ($i = 0; $i -lt $accVals.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 13).Value = $accVals[$i]   # column M = Acc
    $ws.Cells.Item($row, 14).Value = $lossVals[$i]  # column N = Loss
    $ws.Cells.Item($row, 15).Value = $timeVals[$i]  # column O = Time
}

# Restore the active-cell selection recorded in the saved workbook
[void]$ws.Range("R10").Select()
